$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.6606524410359556, 1.655778082260271, 3.537761648806719, 10.19245300693656, 16.0466451790395)
    3 = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 4.358119930609447)
    4 = @(0.6606524410359556, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 2.960089034096801)
    5 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 15.88780690183548)
    6 = @(1.455362044514542, 0.306821227259698, 0.7527432677738641, 10.19245300693656, 12.70737954648466)
    7 = @(0.1190320826869504, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 1.672833113781282)
    8 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E
    $ws.Cells.Item($row, 7).Value = $vals[4]  # G
}
